$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 23
$ws.Range("H4").Value = 23
$ws.Range("E5").Value = 109
$ws.Range("F5").Value = 67
$ws.Range("H5").Value = 67
$ws.Range("F6").Value = 23
$ws.Range("H6").Value = 23
$ws.Range("E10").Value = 372
$ws.Range("F10").Value = 173
$ws.Range("H10").Value = 173
$ws.Range("E11").Value = 248
$ws.Range("F11").Value = 132
$ws.Range("H11").Value = 132
$ws.Range("E12").Value = 363
$ws.Range("F12").Value = 203
$ws.Range("H12").Value = 203
$ws.Range("F13").Value = 54
$ws.Range("H13").Value = 54
$ws.Range("E15").Value = 121
$ws.Range("F15").Value = 48
$ws.Range("H15").Value = 48
$ws.Range("E16").Value = 161
$ws.Range("F16").Value = 79
$ws.Range("H16").Value = 79
$ws.Range("E20").Value = 73
$ws.Range("E21").Value = 118
$ws.Range("E22").Value = 140
$ws.Range("E23").Value = 160
$ws.Range("F23").Value = 71
$ws.Range("H23").Value = 71
$ws.Range("E24").Value = 162
$ws.Range("F24").Value = 83
$ws.Range("H24").Value = 83
$ws.Range("F26").Value = 64
$ws.Range("H26").Value = 64
$ws.Range("E27").Value = 249
$ws.Range("F27").Value = 120
$ws.Range("H27").Value = 120
$ws.Range("E28").Value = 148
$ws.Range("F28").Value = 49
$ws.Range("H28").Value = 49
$ws.Range("F29").Value = 78
$ws.Range("H29").Value = 78
$ws.Range("F30").Value = 91
$ws.Range("H30").Value = 91
$ws.Range("F31").Value = 28
$ws.Range("H31").Value = 28
$ws.Range("E32").Value = 149
$ws.Range("F32").Value = 82
$ws.Range("H32").Value = 82
$ws.Range("F33").Value = 112
$ws.Range("H33").Value = 112
$ws.Range("F34").Value = 98
$ws.Range("H34").Value = 98
$ws.Range("E35").Value = 110
$ws.Range("F35").Value = 69
$ws.Range("H35").Value = 69
$ws.Range("E38").Value = 77
$ws.Range("F39").Value = 72
$ws.Range("H39").Value = 72
$ws.Range("F40").Value = 87
$ws.Range("H40").Value = 87
$ws.Range("E41").Value = 304
$ws.Range("F41").Value = 129
$ws.Range("H41").Value = 129
$ws.Range("E42").Value = 270
$ws.Range("F42").Value = 145
$ws.Range("H42").Value = 145
$ws.Range("E43").Value = 90
$ws.Range("F43").Value = 44
$ws.Range("H43").Value = 44
$ws.Range("E44").Value = 249
$ws.Range("F44").Value = 116
$ws.Range("H44").Value = 116
$ws.Range("E45").Value = 108
$ws.Range("F45").Value = 48
$ws.Range("H45").Value = 48
$ws.Range("E46").Value = 236
$ws.Range("F46").Value = 127
$ws.Range("H46").Value = 127
$ws.Range("E47").Value = 345
$ws.Range("F47").Value = 163
$ws.Range("H47").Value = 163
$ws.Range("E48").Value = 159
$ws.Range("F48").Value = 63
$ws.Range("H48").Value = 63
$ws.Range("E49").Value = 229
$ws.Range("F49").Value = 95
$ws.Range("H49").Value = 95
$ws.Range("E50").Value = 193
$ws.Range("F50").Value = 70
$ws.Range("H50").Value = 70
$ws.Range("F51").Value = 75
$ws.Range("H51").Value = 75
$ws.Range("E52").Value = 22
$ws.Range("F52").Value = 10
$ws.Range("H52").Value = 10
